# Generate Report for Handoff
# Refresh the "d6f573e9-02fc-409c-9f24-201cfa174c29" row's handoff/handback
# timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-08-17 10:08:41"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-17 10:08:37"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-17 10:08:41"
